# Update cryptocurrency price/volume data per latest GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.404.79"
$ws.Range("E2").Value = "  -0.08%  "

# Row 3
$ws.Range("D3").Value = "1.916.23"
$ws.Range("E3").Value = "  +0.76%  "

# Row 4
$ws.Range("E4").Value = "  +0.40%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.89"
$ws.Range("E5").Value = "  -0.27%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  +0.39%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4813"
$ws.Range("E7").Value = "  +0.29%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4069"
$ws.Range("E8").Value = "  +0.07%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08227"
$ws.Range("E9").Value = "  +1.85%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.013"
$ws.Range("E10").Value = "  +0.93%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.32"
$ws.Range("E11").Value = "  +0.03%  "

# Row 12
$ws.Range("D12").Value = "1.925.02"
$ws.Range("E12").Value = "  +2.13%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.064"
$ws.Range("E13").Value = "  +1.93%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.242"
$ws.Range("E14").Value = "  +2.41%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.51"
$ws.Range("E15").Value = "  +1.82%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06836"
$ws.Range("E16").Value = "  +2.17%  "

# Row 17
$ws.Range("E17").Value = "  +0.50%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001037"
$ws.Range("E18").Value = "  +0.65%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.57"
$ws.Range("E19").Value = "  -0.24%  "

# Row 20
$ws.Range("E20").Value = "  +0.48%  "

# Row 21
$ws.Range("D21").Value = "29.418.83"
$ws.Range("E21").Value = "  -0.03%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.649"
$ws.Range("E22").Value = "  +2.07%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.74"
$ws.Range("E23").Value = "  -0.15%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.175"
$ws.Range("E24").Value = "  +0.68%  "

# Row 25
$ws.Range("D25").Value = "2.156.32"
$ws.Range("E25").Value = "  +1.68%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.651"
$ws.Range("E26").Value = "  +9.75%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.50"
$ws.Range("E27").Value = "  +0.43%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.00"
$ws.Range("E28").Value = "  +1.33%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.113"
$ws.Range("E29").Value = "  +0.96%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.52"
$ws.Range("E30").Value = "  +1.90%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.017"
$ws.Range("E31").Value = "  -1.14%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09600"
$ws.Range("E32").Value = "  +0.91%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.656"
$ws.Range("E33").Value = "  +4.82%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.548"
$ws.Range("E34").Value = "  +0.38%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.372"
$ws.Range("E35").Value = "  -1.41%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02280"
$ws.Range("E36").Value = "  +1.15%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06099"
$ws.Range("E37").Value = "  +0.55%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.177"
$ws.Range("E38").Value = "  +0.31%  "

# Row 39
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.065"
$ws.Range("E39").Value = "  +2.43%  "

# Row 40
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5980"
$ws.Range("E40").Value = "  +1.74%  "

# Row 41
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.85"
$ws.Range("E41").Value = "  +6.23%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1845"
$ws.Range("E42").Value = "  +0.05%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.283"
$ws.Range("E43").Value = "  -0.09%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.390"
$ws.Range("E44").Value = "  -1.34%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07594"
$ws.Range("E45").Value = "  -1.92%  "

# Row 46
$ws.Range("E46").Value = "  +0.92%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5591"
$ws.Range("E47").Value = "  +1.25%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.952"
$ws.Range("E48").Value = "  +1.59%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "117.86"
$ws.Range("E49").Value = "  +3.97%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.426"
$ws.Range("E50").Value = "  +3.80%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.21"
$ws.Range("E51").Value = "  +0.12%  "
